# Assignment4_Lobsters.docx edit script
# Commit: "put lobster size data into case format"
#
#  1. select(YEAR, MONTH, DATE, SITE, SIZE, COUNT) -> select(YEAR, SITE, SIZE, COUNT)
#     and add a blank "  " line after the comment that follows it.
#  2. select(YEAR, MONTH, DATE, SITE, TRAPS)       -> select(YEAR, SITE, TRAPS)
#  3. Add a "#4 create graphs ..." comment (plus blank line) before the
#     `abundance_col <-` source chunk.
#  4. Append a brand-new SourceCode paragraph building `lobster_case_format`.

$d = $word.ActiveDocument

function Add-StyledRun {
    # Inserts $text at $pos, tags it with character style $style (if any),
    # and returns the position right after the inserted text.
    param($doc, $pos, $text, $style)
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    $newPos = $pos + $text.Length
    if ($style) {
        $styledRange = $doc.Range($pos, $newPos)
        $styledRange.Style = $style
    }
    return $newPos
}

function Add-LineBreak {
    # Inserts a textWrapping line break (<w:br/>) at $pos; returns new position.
    param($doc, $pos)
    $r = $doc.Range($pos, $pos)
    $r.InsertBreak(6) | Out-Null
    return $pos + 1
}

function Add-RunSequence {
    # Walks an array of @{Kind="text";Style=...;Text=...} / @{Kind="break"}
    # entries starting at $pos; returns the final position.
    param($doc, $pos, $entries)
    foreach ($e in $entries) {
        if ($e.Kind -eq "break") {
            $pos = Add-LineBreak $doc $pos
        } else {
            $pos = Add-StyledRun $doc $pos $e.Text $e.Style
        }
    }
    return $pos
}

# ---------------------------------------------------------------------------
# 1) (YEAR, MONTH, DATE, SITE, SIZE, COUNT)  ->  (YEAR, SITE, SIZE, COUNT)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "(YEAR, MONTH, DATE, SITE, SIZE, COUNT) ", $true, $false, $false, $false,
    $false, $true, 1, $false, "(YEAR, SITE, SIZE, COUNT) ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 1b) Insert a "  " run + extra line break right after the comment that
#     follows the select() above (between the two existing breaks).
# ---------------------------------------------------------------------------
$searchRange = $d.Content.Duplicate
$searchRange.Find.Execute(
    "# Edit Lobster Size data to remove -99999 data from size column and remove transect information"
) | Out-Null
$insPos = $searchRange.End + 1   # right after the first <w:br/> run

$insPos = Add-RunSequence $d $insPos @(
    @{ Kind = "text"; Style = "NormalTok"; Text = "  " }
    @{ Kind = "break" }
)

# ---------------------------------------------------------------------------
# 2) (YEAR, MONTH, DATE, SITE, TRAPS)  ->  (YEAR, SITE, TRAPS)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "(YEAR, MONTH, DATE, SITE, TRAPS)", $true, $false, $false, $false,
    $false, $true, 1, $false, "(YEAR, SITE, TRAPS)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Prepend a "#4 create graphs ..." comment + blank line before
#    "abundance_col <-"
# ---------------------------------------------------------------------------
$abundanceRange = $d.Content.Duplicate
$abundanceRange.Find.Execute("abundance_col <-") | Out-Null
$pos = $abundanceRange.Start

Add-RunSequence $d $pos @(
    @{ Kind = "text"; Style = "CommentTok"; Text = "#4 create graphs of abundance and fishing pressure from 2012-2017" }
    @{ Kind = "break" }
    @{ Kind = "break" }
) | Out-Null

# ---------------------------------------------------------------------------
# 4) New trailing SourceCode paragraph: lobster_case_format
# ---------------------------------------------------------------------------
$newP = $d.Paragraphs.Add()
$newP.Style = "SourceCode"
$pos = $newP.Range.Start

$entries = @(
    @{ Kind = "text"; Style = "CommentTok"; Text = "#put lobster size data into case format." }
    @{ Kind = "break" }
    @{ Kind = "break" }
    @{ Kind = "text"; Style = "NormalTok"; Text = "lobster_case_format <-" }
    @{ Kind = "text"; Style = "StringTok"; Text = " " }
    @{ Kind = "text"; Style = "NormalTok"; Text = "lobster_size_edits " }
    @{ Kind = "text"; Style = "OperatorTok"; Text = "%>%" }
    @{ Kind = "break" }
    @{ Kind = "text"; Style = "StringTok"; Text = "  " }
    @{ Kind = "text"; Style = "KeywordTok"; Text = "filter" }
    @{ Kind = "text"; Style = "NormalTok"; Text = "(YEAR" }
    @{ Kind = "text"; Style = "OperatorTok"; Text = "==" }
    @{ Kind = "text"; Style = "StringTok"; Text = "`"2012`"" }
    @{ Kind = "text"; Style = "NormalTok"; Text = " " }
    @{ Kind = "text"; Style = "OperatorTok"; Text = "|" }
    @{ Kind = "text"; Style = "NormalTok"; Text = "YEAR" }
    @{ Kind = "text"; Style = "OperatorTok"; Text = "==" }
    @{ Kind = "text"; Style = "StringTok"; Text = "`"2017`"" }
    @{ Kind = "text"; Style = "NormalTok"; Text = ") " }
    @{ Kind = "text"; Style = "OperatorTok"; Text = "%>%" }
    @{ Kind = "break" }
    @{ Kind = "text"; Style = "StringTok"; Text = "  " }
    @{ Kind = "text"; Style = "KeywordTok"; Text = "as.data.frame" }
    @{ Kind = "text"; Style = "NormalTok"; Text = "() " }
    @{ Kind = "text"; Style = "OperatorTok"; Text = "%>%" }
    @{ Kind = "break" }
    @{ Kind = "text"; Style = "StringTok"; Text = "  " }
    @{ Kind = "text"; Style = "KeywordTok"; Text = "expand.dft" }
    @{ Kind = "text"; Style = "NormalTok"; Text = "(" }
    @{ Kind = "text"; Style = "DataTypeTok"; Text = "freq=" }
    @{ Kind = "text"; Style = "StringTok"; Text = "`"COUNT`"" }
    @{ Kind = "text"; Style = "NormalTok"; Text = ")" }
)

Add-RunSequence $d $pos $entries | Out-Null

Write-Output "Edit complete."
